# Apply the described edit:
#  - Add F26 = "Kз"
#  - Add F27 = "Kз.ав"
#  - Add K16 = "I35" (text), L16 = K15/(3^(1/2)*35)
#  - Add K17 = "I6"  (text), L17 = K15/(3^(1/2)*6)
#  - Move the view: topLeftCell E8, selection L16
# (new shared strings are written in this exact order so the resulting
#  sharedStrings.xml table lines up with the source table: Kз, Kз.ав, I35, I6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F26").Value = "Kз"
$ws.Range("F27").Value = "Kз.ав"

$ws.Range("K16").Value = "I35"
$ws.Range("L16").Formula = "=K15/(3^(1/2)*35)"

$ws.Range("K17").Value = "I6"
$ws.Range("L17").Formula = "=K15/(3^(1/2)*6)"

$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L16").Select()
